# Rename the worksheet to reflect the functional-testing scope of the test cases.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Name = "Testarea funcțională"

# Move the viewport/selection to the area the author was last working in
# (around the 9th test case block, rows ~104-111).
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 104
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("D106").Select()
